$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Apply replacements from the end of the document towards the start. The
# Find/Replace implementation folds runs that follow a replacement, so
# working backwards keeps the not-yet-visited (earlier) text's run
# structure untouched until it is its own turn to be edited.

Replace-Text "it will be reliable when forecasting the energy consumption." "it will be reliable when forecasting energy consumption."

Replace-Text "Australia’s largest electricity market which ensures Australians have access to reliable" "Australia’s largest electricity market, ensuring Australians have access to reliable"

Replace-Text "as the production cost increases, so does the retail price." "as the production cost increases, so do the retail price."

Replace-Text "Thursday, which has a significant impact on the retail price." "Thursday, which significantly impacts the retail price."

Replace-Text "the organization generates large amount of electricity on Thursday" "the organization generates a large amount of electricity on Thursday"

Replace-Text "AEMO operates the electricity markets by allowing energy related services" "AEMO operates the electricity markets by allowing energy-related services"

Replace-Text "the daytime temperature is typically higher and climate change has a significant impact on the activities Australians engage in." "the daytime temperature is typically higher, and climate change has a significant impact on Australians' activities."

Replace-Text "The price data is reviewed in accordance with National Electricity Rule" "The price data is reviewed following National Electricity Rule"

Replace-Text "Various factors affect the retail price of the electricity and prices become final" "Various factors affect the retail price of the electricity, and prices become final"

Replace-Text "The organization provides affordable electricity and it can be observed" "The organization provides affordable electricity, and it can be observed"

Replace-Text "pandemic has increased the energy consumption" "pandemic have increased the energy consumption"

Replace-Text "From the below graph, it can be seen that during some intervals of the day, the RRP was negative." "From the below graph, it can be seen that the RRP was negative during some intervals of the day."

Replace-Text "RRP value changes on the half-hourly basis" "RRP value changes on a half-hourly basis"

Replace-Text "per megawatt which is set by AEMO" "per megawatt, which AEMO sets"
